$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "RS" column (D) values from 29 to 12 for rows 2 through 8
$ws.Range("D2:D8").Value = 12
